$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0) Remove the pre-existing "_GoBack" bookmark (it currently sits right
#    after the "Shiny app" hyperlink, near the end of the document). It is
#    relocated later in this script, so take it out first while its name is
#    still unambiguous (only one "_GoBack" bookmark exists at this point).
# ---------------------------------------------------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# ---------------------------------------------------------------------------
# 1) "Supplementary Materials: SOR for the manuscript." becomes
#    "Deviation from Sampling Strategy: preregistered deviation from
#    sampling strategy" (the bold label run and the plain-text run are each
#    updated in place so their character formatting is preserved).
# ---------------------------------------------------------------------------
$find1 = $d.Content
$null = $find1.Find.Execute("Supplementary Materials:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraA = $find1.Paragraphs(1)
$rA = $paraA.Range
$null = $rA.Find.Execute("Supplementary Materials:", $true, $false, $false, $false, $false, $true, 1, $false, "Deviation from Sampling Strategy:", 2)
$rA2 = $paraA.Range
$null = $rA2.Find.Execute(" SOR for the manuscript. ", $true, $false, $false, $false, $false, $true, 1, $false, " preregistered deviation from sampling strategy", 2)

# ---------------------------------------------------------------------------
# 2) The old "Deviation from Sampling Strategy: preregistered deviation from
#    sampling strategy" paragraph (now duplicated by step 1 above) is
#    removed entirely, including its paragraph mark.
# ---------------------------------------------------------------------------
$find2 = $d.Content
$null = $find2.Find.Execute("Deviation from Sampling Strategy:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraB = $find2.Paragraphs(1)
$paraB.Range.Delete()

# ---------------------------------------------------------------------------
# 3) The "Supplementary materials for the main article: Additional details
#    ..." bullet gets extra trailing text, with "etc" re-wrapped in
#    proofing-error markers (as Word's spell checker would leave them) and
#    a fresh "_GoBack" bookmark placed at the very end of the paragraph.
# ---------------------------------------------------------------------------
$find3 = $d.Content
$null = $find3.Find.Execute("Supplementary materials for the main article", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraC = $find3.Paragraphs(1)

$paraCXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Supplementary materials for the main article</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Additional details of methods and results, deviations from preregistration, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>; both Supplementary Online Material (Revised; SOM-R) and documentation of changes to the manuscript between Stage 1 and Stage 2.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$null = $paraC.Range.InsertXML($paraCXml)
